$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2:N5").Value = 52.47848103381103
